$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: x input (C3:R3) all bits set to 0 ---
$ws.Range("C3:R3").Value = 0

# --- Row 4 (zx): turn on checkbox ---
$ws.Range("B4").Value = 1

# --- Row 6 (zy): turn on checkbox ---
$ws.Range("B6").Value = 1

# --- Row 7 (ny): turn on checkbox ---
$ws.Range("B7").Value = 1

# --- Row 8 (f): relabel and change formula to add-with-overflow-wrap ---
$ws.Range("A8").Value = "f (1 for add)"
$ws.Range("C8").Formula = '=IF($B$8=1,IF(C7+S7=2,0,C7+S7),IF(C7+S7=2,1,0))'
$ws.Range("D8:R8").Formula = '=IF($B$8=1,IF(D7+T7=2,0,D7+T7),IF(D7+T7=2,1,0))'

# --- Row 9 (no): turn on checkbox and add invert formulas (previously blank) ---
$ws.Range("B9").Value = 1
$ws.Range("C9").Formula = '=IF($B$9=1,IF(C8=1,0,1),C8)'
$ws.Range("D9:R9").Formula = '=IF($B$9=1,IF(D8=1,0,1),D8)'

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 15.0
$ws.Columns("B").ColumnWidth = 9.166666666666666

# --- Selection / view ---
$ws.Range("L10").Select()
